$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142; everything currently at row 142 and below
# shifts down by one row (old row 142 -> new row 143, ..., old row 205 -> new row 206).
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new weekly data point.
$ws.Range("A142").Value = 3
$ws.Range("B142").Value = "Femacal de La Calera"
$ws.Range("C142").Value = "Coquimbo"
$ws.Range("D142").Value = 44510
$ws.Range("E142").Value = 5
$ws.Range("F142").Value = 100112039
$ws.Range("G142").Value = "Ciboulette"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 160
$ws.Range("K142").Value = 1500
$ws.Range("L142").Value = 1500
$ws.Range("M142").Value = 1500
$ws.Range("N142").Value = "`$/docena de atados"
$ws.Range("O142").Value = "Provincia de Quillota"
$ws.Range("P142").Value = 500
$ws.Range("Q142").Value = 3
$ws.Range("R142").Value = "Hortaliza"
